# Disable "smart quotes" autocorrect so literal straight/curly quote
# characters inserted via Find/Replace are preserved exactly as specified
# instead of being auto-converted.
$word.Options.AutoFormatReplaceQuotes = $false

$d = $word.ActiveDocument

function Split-ParagraphWithBreaks($paraIndex, [string[]]$boundaries) {
    # $boundaries is an array of strings of the form "before|after" where a
    # manual line break (<w:br/>) must be inserted between "before" and
    # "after" inside the given paragraph.
    foreach ($b in $boundaries) {
        $parts = $b -split '\|', 2
        $before = $parts[0]
        $after = $parts[1]
        $range = $d.Paragraphs.Item($paraIndex).Range
        $searchText = "$before$after"
        $replaceText = "$before^l$after"
        # Wrap=0 (wdFindStop) keeps the search confined to the supplied
        # Range instead of spilling over into the whole document once the
        # end of the range is reached.
        $found = $range.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 0, $false, $replaceText, 2)
        if (-not $found) {
            throw "Could not find boundary '$searchText' in paragraph $paraIndex"
        }
    }
}

# --- "Programa" (Portuguese) paragraph ---
Split-ParagraphWithBreaks 14 @(
    "RStudio.|2. O que",
    "`"R”.|3. Comandos",
    "básicos.|4. Uso",
    "de dados.|5. Exploração",
    "dos dados.|6. Manipular",
    "utilizando R.|7. Correlação",
    "e regressão.|8. Utilizar",
    "estatísticos.|9. Produzir",
    "ambiente R.|10. Estudos"
)

# --- "Programa" (English) paragraph ---
Split-ParagraphWithBreaks 15 @(
    "software.|2. What",
    "language?|3. Basic",
    "commands.|4. Use",
    "data entry.|5. Preliminary",
    "exploration.|6. Manipulate",
    "using R.|7. Correlation",
    "regression.|8. Use",
    "statistical tests.|9. Produce",
    "environment.|10. Case"
)

# --- "Bibliografia" paragraph ---
Split-ParagraphWithBreaks 19 @(
    "RCommander/|JAMES,",
    "springer, 2013. |MAINDONALD,",
    "Press, 2006. |REYES,",
    "Press, 2017. |TEAM,",
    "(2014).|VENABLES,"
)
